$d = $word.ActiveDocument

$replacements = @(
    @{old="35×79=2765"; new="13×57=741"},
    @{old="45×98=4410"; new="25×80=2000"},
    @{old="43×32=1376"; new="22×20=440"},
    @{old="14×99=1386"; new="65×86=5590"},
    @{old="90×52=4680"; new="97×85=8245"},
    @{old="92×40=3680"; new="77×92=7084"},
    @{old="94×97=9118"; new="32×13=416"},
    @{old="33×19=627";  new="65×21=1365"},
    @{old="11×92=1012"; new="63×22=1386"},
    @{old="30×73=2190"; new="74×94=6956"},
    @{old="27×45=1215"; new="28×23=644"},
    @{old="64×55=3520"; new="24×14=336"},
    @{old="12×73=876";  new="29×50=1450"},
    @{old="69×24=1656"; new="67×63=4221"},
    @{old="41×90=3690"; new="79×62=4898"},
    @{old="38×33=1254"; new="38×17=646"},
    @{old="34×38=1292"; new="62×78=4836"},
    @{old="84×43=3612"; new="74×22=1628"},
    @{old="66×96=6336"; new="24×62=1488"},
    @{old="59×69=4071"; new="38×80=3040"},
    @{old="98×50=4900"; new="13×44=572"},
    @{old="53×40=2120"; new="85×77=6545"},
    @{old="36×67=2412"; new="70×80=5600"},
    @{old="34×99=3366"; new="72×78=5616"},
    @{old="53×30=1590"; new="74×65=4810"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
